$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.080.29"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "1.832.89"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.80"
$ws.Range("E5").Value = "  +0.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6285"
$ws.Range("E6").Value = "  +0.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07477"
$ws.Range("E8").Value = "  -1.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2923"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.10"
$ws.Range("E10").Value = "  +2.04%  "

$ws.Range("E11").Value = "  -0.52%  "

$ws.Range("D12").Value = "1.824.77"
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.975"
$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6686"
$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.64"
$ws.Range("E15").Value = "  -0.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009371"
$ws.Range("E16").Value = "  -5.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.036"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").Value = "29.105.55"
$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("D19").Value = "2.093.49"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.60"
$ws.Range("E20").Value = "  +2.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "222.97"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.136"
$ws.Range("E23").Value = "  -0.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.18"
$ws.Range("E25").Value = "  +0.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1399"
$ws.Range("E26").Value = "  +2.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.500"
$ws.Range("E27").Value = "  +0.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.92"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.501"
$ws.Range("E29").Value = "  +0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05720"
$ws.Range("E30").Value = "  +9.96%  "

$ws.Range("E31").Value = "  +1.84%  "

$ws.Range("E32").Value = "  +0.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.206"
$ws.Range("E33").Value = "  +0.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7491"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.853"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.137"
$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("E37").Value = "  -3.56%  "

$ws.Range("D38").Value = "1.229.82"
$ws.Range("E38").Value = "  -1.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.751"
$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.564"
$ws.Range("E41").Value = "  +3.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8934"
$ws.Range("E42").Value = "  -0.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.09"
$ws.Range("E44").Value = "  +0.62%  "

$ws.Range("D45").Value = "1.985.24"
$ws.Range("E45").Value = "  +0.29%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000125"
$ws.Range("E46").Value = "  -2.00%  "

$ws.Range("B47").Value = "XinFinNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07910"
$ws.Range("E47").Value = "  +16.60%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.56"
$ws.Range("E48").Value = "  +2.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5085"
$ws.Range("E49").Value = "  -0.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4057"
$ws.Range("E50").Value = "  +1.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.046"
$ws.Range("E51").Value = "  +2.54%  "
